# Update the EPEX spot prices workbook with the latest daily data.
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add a new day column (BI) with 13-aug prices ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Mirror the header formatting used by the other day columns (bold,
# bordered, centered) before writing the new header text.
$wsSpot.Range("BH1").Copy()
$wsSpot.Range("BI1").PasteSpecial(-4122)  # xlPasteFormats
$wsSpot.Range("BI1").Value = "13-aug"

$spotValues = @{
    2  = 98.18000000000001
    3  = 90
    4  = 86.89
    5  = 84.70999999999999
    6  = 87.06999999999999
    7  = 92.76000000000001
    8  = 105
    9  = 103.24
    10 = 105.79
    11 = 94.09999999999999
    12 = 90
    13 = 71
    14 = 63.48
    15 = 47.01
    16 = 42.99
    17 = 70.90000000000001
    18 = 84.3
    19 = 96.64
    20 = 115.78
    21 = 150.09
    22 = 155.1
    23 = 153
    24 = 123.21
    25 = 106.27
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Cells.Item($row, 61).Value = $spotValues[$row]
}

# --- Sheet "Gaz": append new row 58 for 2025-08-11 ---
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force text storage (no auto date-conversion) for the new date label,
# then restore the same plain formatting the other date cells carry.
$wsGaz.Range("A58").NumberFormat = "@"
$wsGaz.Range("A58").Value = "2025-08-11"
$wsGaz.Range("A57").Copy()
$wsGaz.Range("A58").PasteSpecial(-4122)  # xlPasteFormats
$wsGaz.Range("B58").Value = 32.4

# --- Sheet "CO2": append new row 58 for 2025-08-11 ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A58").NumberFormat = "@"
$wsCo2.Range("A58").Value = "2025-08-11"
$wsCo2.Range("A57").Copy()
$wsCo2.Range("A58").PasteSpecial(-4122)  # xlPasteFormats
$wsCo2.Range("B58").Value = 71.73999999999999
